# Refresh the "cryptos" price table (Price column D, Volume(1h) column E)
# with the latest scraped values from the GitHub Actions run.
#
# Note: some Price values look numeric (e.g. "309.70", "88.55") but must stay
# literal text (they are inline/shared strings in the sheet, not numbers) so
# that things like a trailing zero are preserved exactly. A leading "'" forces
# Excel to store the literal text instead of silently coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.879.41'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.842.43'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = "'309.70"
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = "'0.4749"
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("E8").Value = '  +2.05%  '
$ws.Range("D9").Value = "'0.07195"
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("D10").Value = "'0.9260"
$ws.Range("E10").Value = '  +3.22%  '
$ws.Range("D11").Value = "'19.62"
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").Value = "'0.07684"
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").Value = '1.882.20'
$ws.Range("E13").Value = '  +2.87%  '
$ws.Range("D14").Value = "'5.299"
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").Value = "'88.55"
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = "'0.000008611"
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = '26.909.94'
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("E21").Value = '  +2.92%  '
$ws.Range("D22").Value = "'5.048"
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("D23").Value = "'10.63"
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = "'152.35"
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("D26").Value = "'18.13"
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").Value = "'1.997"
$ws.Range("E27").Value = '  +1.67%  '
$ws.Range("D28").Value = "'114.10"
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("D29").Value = "'4.921"
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").Value = "'3.307"
$ws.Range("E31").Value = '  +5.56%  '
$ws.Range("D32").Value = "'0.7479"
$ws.Range("E32").Value = '  +2.79%  '
$ws.Range("D33").Value = "'1.169"
$ws.Range("E33").Value = '  +3.98%  '
$ws.Range("D34").Value = "'4.475"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("D36").Value = "'1.089"
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("E37").Value = '  +1.54%  '
$ws.Range("D38").Value = "'0.05258"
$ws.Range("E38").Value = '  +3.11%  '
$ws.Range("D39").Value = "'2.964"
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("D40").Value = "'0.5194"
$ws.Range("E40").Value = '  +3.22%  '
$ws.Range("D41").Value = "'6.958"
$ws.Range("E41").Value = '  +1.74%  '
$ws.Range("D42").Value = "'0.1508"
$ws.Range("E42").Value = '  +1.24%  '
$ws.Range("D43").Value = "'8.193"
$ws.Range("E43").Value = '  +3.09%  '
$ws.Range("D44").Value = "'10.53"
$ws.Range("E44").Value = '  +6.25%  '
$ws.Range("D45").Value = "'0.4723"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = "'101.37"
$ws.Range("E47").Value = '  +3.47%  '
$ws.Range("E48").Value = '  +2.92%  '
$ws.Range("D49").Value = "'65.32"
$ws.Range("E49").Value = '  +2.74%  '
$ws.Range("D50").Value = "'0.06019"
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").Value = "'0.8860"
$ws.Range("E51").Value = '  +4.45%  '
